$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'',
                                                     random_state=42))),
                (''model'',
                 BaggingClassifier(estimator=SVC(C=0.0001, kernel=''linear'',
                                                 random_state=42),
                                   n_estimators=50, random_state=42))])'
$ws.Range('B2').Value = 0.7499999999999999
$ws.Range('C2').Value = '{''selector'': SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'', random_state=42)), ''scaler'': MinMaxScaler(), ''model__n_estimators'': 50, ''model__estimator__kernel'': ''linear'', ''model__estimator__class_weight'': None, ''model__estimator__C'': 0.0001}'
$ws.Range('D2').Value = 0.6758066402468806
$ws.Range('E2').Value = 0.5506162892662892
$ws.Range('F2').Value = 0.8
$ws.Range('G2').Value = 0.6299342502011319
$ws.Range('H2').Value = 0.50355873015873
$ws.Range('I2').Value = 0.6666666666666666
$ws.Range('J2').Value = 0.7798085106382977
$ws.Range('K2').Value = 0.6519999999999999
$ws.Range('L2').Value = 1
$ws.Range('M2').Value = '[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1]'
$ws.Range('N2').Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'
$ws.Range('O2').Value = 42

$ws.Range('A3').Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                (''model'',
                 BaggingClassifier(estimator=SVC(C=0.0001, kernel=''poly'',
                                                 random_state=42),
                                   random_state=42))])'
$ws.Range('B3').Value = 0.7499999999999999
$ws.Range('C3').Value = '{''selector'': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), ''scaler'': MinMaxScaler(), ''model__n_estimators'': 10, ''model__estimator__kernel'': ''poly'', ''model__estimator__class_weight'': None, ''model__estimator__C'': 0.0001}'
$ws.Range('D3').Value = 0.6704001125603487
$ws.Range('E3').Value = 0.5480335941835941
$ws.Range('F3').Value = 0.8
$ws.Range('G3').Value = 0.6372086088800542
$ws.Range('H3').Value = 0.5459833333333334
$ws.Range('I3').Value = 0.6666666666666666
$ws.Range('J3').Value = 0.7689574468085106
$ws.Range('K3').Value = 0.6211666666666665
$ws.Range('L3').Value = 1
$ws.Range('M3').Value = '[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0]'
$ws.Range('N3').Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'
$ws.Range('O3').Value = 69

$ws.Range('A4').Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faaa3ed90d0>),
                (''model'',
                 BaggingClassifier(estimator=SVC(C=5, kernel=''linear'',
                                                 random_state=42),
                                   n_estimators=50, random_state=42))])'
$ws.Range('B4').Value = 0.6933333333333332
$ws.Range('C4').Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa9f72b0d0>, ''scaler'': MinMaxScaler(), ''model__n_estimators'': 50, ''model__estimator__kernel'': ''linear'', ''model__estimator__class_weight'': None, ''model__estimator__C'': 5}'
$ws.Range('D4').Value = 0.6615637762024856
$ws.Range('E4').Value = 0.5197645576645578
$ws.Range('F4').Value = 0.7058823529411765
$ws.Range('G4').Value = 0.6219181923063317
$ws.Range('H4').Value = 0.4695087301587302
$ws.Range('I4').Value = 0.8
$ws.Range('J4').Value = 0.7780444444444445
$ws.Range('K4').Value = 0.6630000000000001
$ws.Range('L4').Value = 0.631578947368421
$ws.Range('M4').Value = '[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1]'
$ws.Range('N4').Value = '[1 1 1 0 1 1 0 1 1 0 1 0 0 1 1 1 0 1 0 0 1 0 1 1]'
$ws.Range('O4').Value = 23

$ws.Range('A5').Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9f7d9df0>),
                (''model'',
                 BaggingClassifier(estimator=SVC(C=0.0001, kernel=''linear'',
                                                 random_state=42),
                                   random_state=42))])'
$ws.Range('B5').Value = 0.7499999999999999
$ws.Range('C5').Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa9f8c4be0>, ''scaler'': MinMaxScaler(), ''model__n_estimators'': 10, ''model__estimator__kernel'': ''linear'', ''model__estimator__class_weight'': None, ''model__estimator__C'': 0.0001}'
$ws.Range('D5').Value = 0.6137188472413121
$ws.Range('E5').Value = 0.5140572316572316
$ws.Range('F5').Value = 0.7368421052631579
$ws.Range('G5').Value = 0.5806319370304969
$ws.Range('H5').Value = 0.4749615079365078
$ws.Range('I5').Value = 0.5833333333333334
$ws.Range('J5').Value = 0.7009795918367346
$ws.Range('K5').Value = 0.603
$ws.Range('L5').Value = 1
$ws.Range('M5').Value = '[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1]'
$ws.Range('N5').Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'
$ws.Range('O5').Value = 99

$ws.Range('A6').Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9f53ef70>),
                (''model'',
                 BaggingClassifier(estimator=SVC(C=0.0001, kernel=''poly'',
                                                 random_state=42),
                                   n_estimators=5, random_state=42))])'
$ws.Range('B6').Value = 0.7499999999999999
$ws.Range('C6').Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa9f7eda00>, ''scaler'': MinMaxScaler(), ''model__n_estimators'': 5, ''model__estimator__kernel'': ''poly'', ''model__estimator__class_weight'': None, ''model__estimator__C'': 0.0001}'
$ws.Range('D6').Value = 0.7050542655399822
$ws.Range('E6').Value = 0.5687642524142524
$ws.Range('F6').Value = 0.6285714285714286
$ws.Range('G6').Value = 0.6533995400080819
$ws.Range('H6').Value = 0.519706746031746
$ws.Range('I6').Value = 0.4583333333333333
$ws.Range('J6').Value = 0.7939999999999998
$ws.Range('K6').Value = 0.6699999999999998
$ws.Range('L6').Value = 1
$ws.Range('M6').Value = '[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1]'
$ws.Range('N6').Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'
$ws.Range('O6').Value = 89
